$wb = $excel.ActiveWorkbook

# --- Step 1: swap row 6 and row 7 cell content (columns A-D) on all three sheets ---

$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A6").Value = "cfe47a8c-7699-40ac-bc73-1067185b44ec.md"
$ws.Range("B6").Value = "In Translation"
$ws.Range("C6").Value = "In Translation"
$ws.Range("A7").Value = "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.md"
$ws.Range("B7").Value = "Ready for handoff"
$ws.Range("C7").Value = "Ready for handoff"

$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A6").Value = "cfe47a8c-7699-40ac-bc73-1067185b44ec.md"
$ws.Range("B6").Value = "In Translation"
$ws.Range("C6").Value = "cfe47a8c-7699-40ac-bc73-1067185b44ec.6f50fd140f48297a5e373e207c810326eade4368.zh-cn.xlf"
$ws.Range("D6").Value = "2016-03-01 07:20:30"
$ws.Range("A7").Value = "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.md"
$ws.Range("B7").Value = "Ready for handoff"
$ws.Range("C7").Value = "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.8d3edc7ad4ee28f11f642b67b41601399bece039.zh-cn.xlf"
$ws.Range("D7").Value = "2016-03-01 07:16:36"

$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A6").Value = "cfe47a8c-7699-40ac-bc73-1067185b44ec.md"
$ws.Range("B6").Value = "In Translation"
$ws.Range("C6").Value = "cfe47a8c-7699-40ac-bc73-1067185b44ec.6f50fd140f48297a5e373e207c810326eade4368.de-de.xlf"
$ws.Range("D6").Value = "2016-03-01 07:20:39"
$ws.Range("A7").Value = "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.md"
$ws.Range("B7").Value = "Ready for handoff"
$ws.Range("C7").Value = "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.8d3edc7ad4ee28f11f642b67b41601399bece039.de-de.xlf"
$ws.Range("D7").Value = "2016-03-01 07:16:45"

# --- Step 2: rebuild hyperlinks on all three sheets so that display text follows the new row content ---
# (target URLs / r:id <-> cell ordering is preserved exactly as before the edit; only the
#  display text for the affected cells changes, matching the source diff.)

$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a755bd49dd560c91f302d985a21d1a48a5e0423f/e2e/57e180a2-d839-41a5-8afa-5f903f67d079.md", [Type]::Missing, [Type]::Missing, "57e180a2-d839-41a5-8afa-5f903f67d079.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/663844c1bff88e2ec68a2b5148825bca3ac4c8fd/e2e/6dab1213-1c80-4430-aa58-5f19bf960db5.md", [Type]::Missing, [Type]::Missing, "6dab1213-1c80-4430-aa58-5f19bf960db5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/dbbf376580158e3d229e52ed45936bb16da9dc7d/e2e/74c8d44a-3249-4b3f-afae-70d8bbeca7af.md", [Type]::Missing, [Type]::Missing, "74c8d44a-3249-4b3f-afae-70d8bbeca7af.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/a755bd49dd560c91f302d985a21d1a48a5e0423f/e2e/8d89c7c6-79dd-4651-b6fd-40556e03f5ae.md", [Type]::Missing, [Type]::Missing, "8d89c7c6-79dd-4651-b6fd-40556e03f5ae.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/e0042f5171965ca417079935a01703fdb286632e/e2e/93168050-f6e4-4fe2-8c5c-d15dc8b9b447.md", [Type]::Missing, [Type]::Missing, "cfe47a8c-7699-40ac-bc73-1067185b44ec.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/8ccfa16f51f341fffc48135d0076471a0974c417/e2e/cfe47a8c-7699-40ac-bc73-1067185b44ec.md", [Type]::Missing, [Type]::Missing, "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/8ccfa16f51f341fffc48135d0076471a0974c417/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a755bd49dd560c91f302d985a21d1a48a5e0423f/e2e/57e180a2-d839-41a5-8afa-5f903f67d079.md", [Type]::Missing, [Type]::Missing, "57e180a2-d839-41a5-8afa-5f903f67d079.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c758b2455cecc080017eb4e4c04e5917ad9d2d0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/57e180a2-d839-41a5-8afa-5f903f67d079.6b354ecdb366751c2910c4a9616cf0d0e43510f1.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "57e180a2-d839-41a5-8afa-5f903f67d079.6b354ecdb366751c2910c4a9616cf0d0e43510f1.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/663844c1bff88e2ec68a2b5148825bca3ac4c8fd/e2e/6dab1213-1c80-4430-aa58-5f19bf960db5.md", [Type]::Missing, [Type]::Missing, "6dab1213-1c80-4430-aa58-5f19bf960db5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3316acb1b14bbec0c562f4949962b0d6bb77b77c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/6dab1213-1c80-4430-aa58-5f19bf960db5.dae169b2993541b6ad931e636686067f10d13881.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "6dab1213-1c80-4430-aa58-5f19bf960db5.dae169b2993541b6ad931e636686067f10d13881.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/1df4ad841938ddb8088af986fc426fcf9d54add1/e2e/6dab1213-1c80-4430-aa58-5f19bf960db5.md", [Type]::Missing, [Type]::Missing, "6dab1213-1c80-4430-aa58-5f19bf960db5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e35cdda26b836624e8bf8cb2543ad022b27238e5/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/6dab1213-1c80-4430-aa58-5f19bf960db5.dae169b2993541b6ad931e636686067f10d13881.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "6dab1213-1c80-4430-aa58-5f19bf960db5.dae169b2993541b6ad931e636686067f10d13881.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/dbbf376580158e3d229e52ed45936bb16da9dc7d/e2e/74c8d44a-3249-4b3f-afae-70d8bbeca7af.md", [Type]::Missing, [Type]::Missing, "74c8d44a-3249-4b3f-afae-70d8bbeca7af.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1ea6f006c0bcf2aa1f0ecdb93db861a7ce890540/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/74c8d44a-3249-4b3f-afae-70d8bbeca7af.ec448cd7077d1a98e4f8214315a07b5cbf97e2de.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "74c8d44a-3249-4b3f-afae-70d8bbeca7af.ec448cd7077d1a98e4f8214315a07b5cbf97e2de.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/a755bd49dd560c91f302d985a21d1a48a5e0423f/e2e/8d89c7c6-79dd-4651-b6fd-40556e03f5ae.md", [Type]::Missing, [Type]::Missing, "8d89c7c6-79dd-4651-b6fd-40556e03f5ae.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5c758b2455cecc080017eb4e4c04e5917ad9d2d0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/8d89c7c6-79dd-4651-b6fd-40556e03f5ae.7934dda4d3b5626aa2a4c803e3a7985cf7c9a47d.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "8d89c7c6-79dd-4651-b6fd-40556e03f5ae.7934dda4d3b5626aa2a4c803e3a7985cf7c9a47d.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/e0042f5171965ca417079935a01703fdb286632e/e2e/93168050-f6e4-4fe2-8c5c-d15dc8b9b447.md", [Type]::Missing, [Type]::Missing, "cfe47a8c-7699-40ac-bc73-1067185b44ec.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9ef253d1d63e4477f87bf13ef2b95be806118e36/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/93168050-f6e4-4fe2-8c5c-d15dc8b9b447.8d3edc7ad4ee28f11f642b67b41601399bece039.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "cfe47a8c-7699-40ac-bc73-1067185b44ec.6f50fd140f48297a5e373e207c810326eade4368.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/8ccfa16f51f341fffc48135d0076471a0974c417/e2e/cfe47a8c-7699-40ac-bc73-1067185b44ec.md", [Type]::Missing, [Type]::Missing, "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/12f7b2feb07d95d98954d76e1c6b17a63312e56b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/cfe47a8c-7699-40ac-bc73-1067185b44ec.6f50fd140f48297a5e373e207c810326eade4368.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.8d3edc7ad4ee28f11f642b67b41601399bece039.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/8ccfa16f51f341fffc48135d0076471a0974c417/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null

$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/a755bd49dd560c91f302d985a21d1a48a5e0423f/e2e/57e180a2-d839-41a5-8afa-5f903f67d079.md", [Type]::Missing, [Type]::Missing, "57e180a2-d839-41a5-8afa-5f903f67d079.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/193f1815819a5891391c34cfc7e345bbb874b3d7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/57e180a2-d839-41a5-8afa-5f903f67d079.6b354ecdb366751c2910c4a9616cf0d0e43510f1.de-de.xlf", [Type]::Missing, [Type]::Missing, "57e180a2-d839-41a5-8afa-5f903f67d079.6b354ecdb366751c2910c4a9616cf0d0e43510f1.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/663844c1bff88e2ec68a2b5148825bca3ac4c8fd/e2e/6dab1213-1c80-4430-aa58-5f19bf960db5.md", [Type]::Missing, [Type]::Missing, "6dab1213-1c80-4430-aa58-5f19bf960db5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f168ec5896f70fb58ec32f993dd6d878aaa9597c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/6dab1213-1c80-4430-aa58-5f19bf960db5.dae169b2993541b6ad931e636686067f10d13881.de-de.xlf", [Type]::Missing, [Type]::Missing, "6dab1213-1c80-4430-aa58-5f19bf960db5.dae169b2993541b6ad931e636686067f10d13881.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/985ecf70cab566a695048b39a7b80246fc65fd4a/e2e/6dab1213-1c80-4430-aa58-5f19bf960db5.md", [Type]::Missing, [Type]::Missing, "6dab1213-1c80-4430-aa58-5f19bf960db5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e4ca427b487c63a0ef13722fd9764cf77d6181f8/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/6dab1213-1c80-4430-aa58-5f19bf960db5.dae169b2993541b6ad931e636686067f10d13881.de-de.xlf", [Type]::Missing, [Type]::Missing, "6dab1213-1c80-4430-aa58-5f19bf960db5.dae169b2993541b6ad931e636686067f10d13881.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/dbbf376580158e3d229e52ed45936bb16da9dc7d/e2e/74c8d44a-3249-4b3f-afae-70d8bbeca7af.md", [Type]::Missing, [Type]::Missing, "74c8d44a-3249-4b3f-afae-70d8bbeca7af.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1e9bbec3be59d100fd87376f173b1ec253a832cb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/74c8d44a-3249-4b3f-afae-70d8bbeca7af.ec448cd7077d1a98e4f8214315a07b5cbf97e2de.de-de.xlf", [Type]::Missing, [Type]::Missing, "74c8d44a-3249-4b3f-afae-70d8bbeca7af.ec448cd7077d1a98e4f8214315a07b5cbf97e2de.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/a755bd49dd560c91f302d985a21d1a48a5e0423f/e2e/8d89c7c6-79dd-4651-b6fd-40556e03f5ae.md", [Type]::Missing, [Type]::Missing, "8d89c7c6-79dd-4651-b6fd-40556e03f5ae.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/193f1815819a5891391c34cfc7e345bbb874b3d7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/8d89c7c6-79dd-4651-b6fd-40556e03f5ae.7934dda4d3b5626aa2a4c803e3a7985cf7c9a47d.de-de.xlf", [Type]::Missing, [Type]::Missing, "8d89c7c6-79dd-4651-b6fd-40556e03f5ae.7934dda4d3b5626aa2a4c803e3a7985cf7c9a47d.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/e0042f5171965ca417079935a01703fdb286632e/e2e/93168050-f6e4-4fe2-8c5c-d15dc8b9b447.md", [Type]::Missing, [Type]::Missing, "cfe47a8c-7699-40ac-bc73-1067185b44ec.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/df1e00761f44ce7fe18699341769fff09826761d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/93168050-f6e4-4fe2-8c5c-d15dc8b9b447.8d3edc7ad4ee28f11f642b67b41601399bece039.de-de.xlf", [Type]::Missing, [Type]::Missing, "cfe47a8c-7699-40ac-bc73-1067185b44ec.6f50fd140f48297a5e373e207c810326eade4368.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/8ccfa16f51f341fffc48135d0076471a0974c417/e2e/cfe47a8c-7699-40ac-bc73-1067185b44ec.md", [Type]::Missing, [Type]::Missing, "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b70ee1dfc4579ab1dfdcf123e2f3f87e1665225/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/cfe47a8c-7699-40ac-bc73-1067185b44ec.6f50fd140f48297a5e373e207c810326eade4368.de-de.xlf", [Type]::Missing, [Type]::Missing, "93168050-f6e4-4fe2-8c5c-d15dc8b9b447.8d3edc7ad4ee28f11f642b67b41601399bece039.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "https://github.com/OpenLocalizationTest/oltest/blob/8ccfa16f51f341fffc48135d0076471a0974c417/.localization-config", [Type]::Missing, [Type]::Missing, ".localization-config") | Out-Null
